$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    # Assign via Formula with a leading apostrophe so Excel stores the
    # literal text verbatim (no numeric/date auto-conversion), then reset
    # the style back to Normal so no stray quote-prefix / text-format
    # style sticks to the cell (keeps the original default styling).
    $ws.Range($addr).Formula = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue $ws 'D2' '36.825.37'
Set-TextValue $ws 'E2' '  -0.80%  '
Set-TextValue $ws 'D3' '2.098.35'
Set-TextValue $ws 'E3' '  +2.48%  '
Set-TextValue $ws 'E4' '  -0.15%  '
Set-TextValue $ws 'D5' '245.79'
Set-TextValue $ws 'E5' '  -0.80%  '
Set-TextValue $ws 'E6' '  -1.14%  '
Set-TextValue $ws 'E7' '  -0.02%  '
Set-TextValue $ws 'D8' '55.07'
Set-TextValue $ws 'E8' '  -1.52%  '
Set-TextValue $ws 'D9' '59.18'
Set-TextValue $ws 'E9' '  -0.97%  '
Set-TextValue $ws 'E10' '  -3.00%  '
Set-TextValue $ws 'D11' '0.0766'
Set-TextValue $ws 'E11' '  -1.58%  '
Set-TextValue $ws 'E12' '  +1.32%  '
Set-TextValue $ws 'D13' '0.921'
Set-TextValue $ws 'E13' '  +4.98%  '
Set-TextValue $ws 'E14' '  -6.58%  '
Set-TextValue $ws 'D15' '2.405.99'
Set-TextValue $ws 'E15' '  +2.50%  '
Set-TextValue $ws 'E16' '  -2.79%  '
Set-TextValue $ws 'D17' '2.098.12'
Set-TextValue $ws 'E17' '  +2.36%  '
Set-TextValue $ws 'D18' '36.804.19'
Set-TextValue $ws 'D19' '17.25'
Set-TextValue $ws 'E19' '  -6.23%  '
Set-TextValue $ws 'E20' '  -2.10%  '
Set-TextValue $ws 'D21' '0.0₃0883'
Set-TextValue $ws 'E21' '  -0.87%  '
Set-TextValue $ws 'E22' '  +1.93%  '
Set-TextValue $ws 'D23' '239.18'
Set-TextValue $ws 'E23' '  +1.17%  '
Set-TextValue $ws 'E25' '  -2.70%  '
Set-TextValue $ws 'D26' '9.76'
Set-TextValue $ws 'E26' '  +2.80%  '
Set-TextValue $ws 'D27' '2.17'
Set-TextValue $ws 'E27' '  +0.29%  '
Set-TextValue $ws 'D28' '167.32'
Set-TextValue $ws 'E28' '  -1.15%  '
Set-TextValue $ws 'D29' '20.97'
Set-TextValue $ws 'E29' '  +4.70%  '
Set-TextValue $ws 'E30' '  -1.02%  '
Set-TextValue $ws 'D31' '5.24'
Set-TextValue $ws 'E31' '  +8.42%  '
Set-TextValue $ws 'E32' '  +2.87%  '
Set-TextValue $ws 'D33' '4.73'
Set-TextValue $ws 'E33' '  +5.80%  '
Set-TextValue $ws 'E34' '  -1.21%  '
Set-TextValue $ws 'E35' '  +9.67%  '
Set-TextValue $ws 'E36' '  +0.02%  '
Set-TextValue $ws 'E37' '  +3.65%  '
Set-TextValue $ws 'D38' '0.0821'
Set-TextValue $ws 'E38' '  -7.35%  '
Set-TextValue $ws 'E39' '  -4.36%  '
Set-TextValue $ws 'E40' '  +1.37%  '
Set-TextValue $ws 'B41' 'THORChain'
Set-TextValue $ws 'C41' 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws 'D41' '4.92'
Set-TextValue $ws 'E41' '  -6.53%  '
Set-TextValue $ws 'B42' 'VeChain'
Set-TextValue $ws 'C42' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D42' '0.0221'
Set-TextValue $ws 'E42' '  -0.39%  '
Set-TextValue $ws 'D43' '0.0956'
Set-TextValue $ws 'E43' '  -1.74%  '
Set-TextValue $ws 'D44' '96.71'
Set-TextValue $ws 'E44' '  +1.45%  '
Set-TextValue $ws 'E45' '  -10.01%  '
Set-TextValue $ws 'D46' '1.416.81'
Set-TextValue $ws 'E46' '  +11.98%  '
Set-TextValue $ws 'E47' '  -6.31%  '
Set-TextValue $ws 'D48' '7.60'
Set-TextValue $ws 'E48' '  +12.42%  '
Set-TextValue $ws 'D49' '2.48'
Set-TextValue $ws 'E49' '  +2.37%  '
Set-TextValue $ws 'E50' '  +2.22%  '
Set-TextValue $ws 'D51' '2.292.22'
